# Mark every CRUD test result as Passed (TRUE) on the "Test Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Select the sheet and set the active cell / selection as in the edited file
$ws.Activate()
$ws.Range("H19").Select()

# Data rows 2-24 span columns B (Create) through E (Delete)
$ws.Range("B2:E24").Value = $true
